$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Wnt10b -> Fzd7 -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.377081
$ws.Range("H2").Value = 10.131243
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.418393
$ws.Range("N2").Value = 7.255179
$ws.Range("O2").Value = 0.0919828589765645
$ws.Range("P2").Value = 0.0919828589765645
$ws.Range("Q2").Value = 8.167109050833
$ws.Range("R2").Value = 73.50398145749699
$ws.Range("S2").Value = 0.0919828589765645
$ws.Range("T2").Value = 0.0919828589765645

# Row 3 (FAPs -> Wnt10b -> Fzd7 -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.377081
$ws.Range("H3").Value = 10.131243
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.11799233333333
$ws.Range("N3").Value = 30.353977
$ws.Range("O3").Value = 0.3848348311969811
$ws.Range("P3").Value = 0.3848348311969811
$ws.Range("Q3").Value = 34.16927966704566
$ws.Range("R3").Value = 307.523517003411
$ws.Range("S3").Value = 0.3848348311969811
$ws.Range("T3").Value = 0.3848348311969811

# Row 4 (FAPs -> Wnt10b -> Fzd7 -> sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.377081
$ws.Range("H4").Value = 10.131243
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.75539366666667
$ws.Range("N4").Value = 41.266181
$ws.Range("O4").Value = 0.5231823098264544
$ws.Range("P4").Value = 0.5231823098264544
$ws.Range("Q4").Value = 46.45307859922034
$ws.Range("R4").Value = 418.077707392983
$ws.Range("S4").Value = 0.5231823098264544
$ws.Range("T4").Value = 0.5231823098264544
